$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab20")

# --- Footnote ("resource-rich") reassignment: South Sudan gains the asterisk, Nigeria loses it ---
$ws.Range("B34").Value = "South Sudan*"
$ws.Range("B57").Value = "Nigeria"

# Carry over the "resource-rich" row shading: copy the format of an already-shaded
# data row onto South Sudan's row, and the format of a plain (unshaded) data row
# onto Nigeria's row, without touching the underlying cell values.
$ws.Range("B17:P17").Copy()
$ws.Range("B34:P34").PasteSpecial(-4122)

$ws.Range("B5:P5").Copy()
$ws.Range("B57:P57").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Updated aggregate statistics (recomputed downstream of the resource-rich reclassification) ---
$ws.Range("C69").Value = [double]"1.51327883955708"
$ws.Range("D69").Value = [double]"0.72656260188280997"
$ws.Range("E69").Value = [double]"4.9276778698719097"
$ws.Range("F69").Value = [double]"2.5843251887982501"
$ws.Range("G69").Value = [double]"9.7518445001100496"
$ws.Range("H69").Value = [double]"1.52041912790205"
$ws.Range("I69").Value = [double]"6.5941195910888197"
$ws.Range("J69").Value = [double]"4669.4882442923999"
$ws.Range("K69").Value = [double]"1494.32168469509"
$ws.Range("L69").Value = [double]"15774.18"
$ws.Range("M69").Value = [double]"8428.3237526352696"
$ws.Range("N69").Value = [double]"30366.3136816228"
$ws.Range("O69").Value = [double]"3453.3774651127801"
$ws.Range("P69").Value = [double]"15653.7539360468"
$ws.Range("C77").Value = [double]"2.6965275621193698"
$ws.Range("E77").Value = [double]"0.17671136026486001"
$ws.Range("G77").Value = [double]"5.7080283888388097"
$ws.Range("J77").Value = [double]"88910.267275200007"
$ws.Range("L77").Value = [double]"4627.08"
$ws.Range("N77").Value = [double]"188961.57092620499"
$ws.Range("C80").Value = [double]"0.28710519427525"
$ws.Range("D80").Value = [double]"6.1921681134500003E-6"
$ws.Range("E80").Value = [double]"0.81731710295307003"
$ws.Range("F80").Value = [double]"0.98684404710936002"
$ws.Range("G80").Value = [double]"2.09127253650579"
$ws.Range("H80").Value = [double]"-0.9107816051216"
$ws.Range("I80").Value = [double]"2.4380471240906898"
$ws.Range("J80").Value = [double]"3305.7262295046999"
$ws.Range("K80").Value = [double]"6.8871679038999996E-3"
$ws.Range("L80").Value = [double]"4499.6400000000003"
$ws.Range("M80").Value = [double]"3502.53154962177"
$ws.Range("N80").Value = [double]"11307.9046662944"
$ws.Range("O80").Value = [double]"-2326.4074425386002"
$ws.Range("P80").Value = [double]"7126.9084449928696"
$ws.Range("C82").Value = [double]"2.9842536520177401"
$ws.Range("D82").Value = [double]"-8.3234868015199998E-2"
$ws.Range("E82").Value = [double]"2.7218435374577301"
$ws.Range("F82").Value = [double]"4.1294425424134502"
$ws.Range("G82").Value = [double]"9.75230486387375"
$ws.Range("H82").Value = [double]"3.3688839329938398"
$ws.Range("I82").Value = [double]"2.5872996924920799"
$ws.Range("J82").Value = [double]"79684.812172964594"
$ws.Range("K82").Value = [double]"-9670.2689766842996"
$ws.Range("L82").Value = [double]"60512.38"
$ws.Range("M82").Value = [double]"92655.698919267597"
$ws.Range("N82").Value = [double]"223182.622115548"
$ws.Range("O82").Value = [double]"80198.146930546398"
$ws.Range("P82").Value = [double]"51430.961600495299"
$ws.Range("C84").Value = [double]"3.6374287074057698"
$ws.Range("D84").Value = [double]"0.63438055125051995"
$ws.Range("E84").Value = [double]"7.84991411932715"
$ws.Range("F84").Value = [double]"2.9522167416999801"
$ws.Range("G84").Value = [double]"15.0739401196834"
$ws.Range("H84").Value = [double]"3.7427216219185699"
$ws.Range("I84").Value = [double]"4.4585777752550397"
$ws.Range("J84").Value = [double]"16529.2635764636"
$ws.Range("K84").Value = [double]"1447.5282172213399"
$ws.Range("L84").Value = [double]"34562.080000000002"
$ws.Range("M84").Value = [double]"13370.910330721599"
$ws.Range("N84").Value = [double]"65909.7821244066"
$ws.Range("O84").Value = [double]"13098.0396164182"
$ws.Range("P84").Value = [double]"14861.506007577"
$ws.Range("C86").Value = [double]"1.1351810124637201"
$ws.Range("D86").Value = [double]"1.1458862310186899"
$ws.Range("E86").Value = [double]"1.55850680756045"
$ws.Range("F86").Value = [double]"4.6762991982203301"
$ws.Range("G86").Value = [double]"8.5158732492631994"
$ws.Range("H86").Value = [double]"0.82623000200950003"
$ws.Range("I86").Value = [double]"2.4128981373432699"
$ws.Range("J86").Value = [double]"22569.724626633801"
$ws.Range("K86").Value = [double]"15436.0891897576"
$ws.Range("L86").Value = [double]"28383.61"
$ws.Range("M86").Value = [double]"81460.557384713393"
$ws.Range("N86").Value = [double]"147849.98120110499"
$ws.Range("O86").Value = [double]"13508.955185029799"
$ws.Range("P86").Value = [double]"38246.988668525097"
$ws.Range("C87").Value = [double]"1.6814271499307301"
$ws.Range("D87").Value = [double]"0.31810861600083001"
$ws.Range("E87").Value = [double]"0.39795917852585999"
$ws.Range("F87").Value = [double]"4.2584117087990299"
$ws.Range("G87").Value = [double]"6.6559066532564497"
$ws.Range("H87").Value = [double]"1.8098163404808001"
$ws.Range("I87").Value = [double]"3.1174084655052599"
$ws.Range("J87").Value = [double]"119308.99195780601"
$ws.Range("K87").Value = [double]"19228.405080276301"
$ws.Range("L87").Value = [double]"30429.82"
$ws.Range("M87").Value = [double]"287209.29664095299"
$ws.Range("N87").Value = [double]"456176.51367903501"
$ws.Range("O87").Value = [double]"118624.982667713"
$ws.Range("P87").Value = [double]"189590.902264845"
$ws.Range("C89").Value = [double]"1.5688842946884201"
$ws.Range("D89").Value = [double]"0.74543296815672999"
$ws.Range("E89").Value = [double]"9.4438509363589995E-2"
$ws.Range("F89").Value = [double]"0.89609586035772004"
$ws.Range("G89").Value = [double]"3.3048516325664599"
$ws.Range("H89").Value = [double]"2.14108620481686"
$ws.Range("I89").Value = [double]"1.3996751433230299"
$ws.Range("J89").Value = [double]"372737.11317193299"
$ws.Range("K89").Value = [double]"211924.98088836999"
$ws.Range("L89").Value = [double]"16857.07"
$ws.Range("M89").Value = [double]"208244.05450639699"
$ws.Range("N89").Value = [double]"809763.2185667"
$ws.Range("O89").Value = [double]"553900.21728234598"
$ws.Range("P89").Value = [double]"325092.83516322199"
$ws.Range("C90").Value = [double]"1.56460121754489"
$ws.Range("D90").Value = [double]"4.7374617042655496"
$ws.Range("E90").Value = [double]"0.16583510043134"
$ws.Range("F90").Value = [double]"0.36709999804824001"
$ws.Range("G90").Value = [double]"6.8349980202900298"
$ws.Range("H90").Value = [double]"1.8417994785607501"
$ws.Range("I90").Value = [double]"5.1862699685415397"
$ws.Range("J90").Value = [double]"851070.37732905895"
$ws.Range("K90").Value = [double]"2911844.9666242399"
$ws.Range("L90").Value = [double]"129.68"
$ws.Range("M90").Value = [double]"181987.06973772901"
$ws.Range("N90").Value = [double]"3945032.0936910301"
$ws.Range("O90").Value = [double]"955952.50525997696"
$ws.Range("P90").Value = [double]"2986255.22516564"
$ws.Range("C94").Value = [double]"14.9835775325241"
$ws.Range("G94").Value = [double]"24.602251876796998"
$ws.Range("J94").Value = [double]"106213.991503212"
$ws.Range("N94").Value = [double]"152351.14952677899"
$ws.Range("E95").Value = [double]"5.9524476355772897"
$ws.Range("C97").Value = [double]"1.8190956117673101"
$ws.Range("D97").Value = [double]"0.98331820119424995"
$ws.Range("E97").Value = [double]"3.5806269247697902"
$ws.Range("F97").Value = [double]"3.0232511862690798"
$ws.Range("G97").Value = [double]"9.4062919240004295"
$ws.Range("H97").Value = [double]"1.44670933907934"
$ws.Range("I97").Value = [double]"3.6198534262834099"
$ws.Range("J97").Value = [double]"25499.395600239201"
$ws.Range("K97").Value = [double]"8951.5245118268103"
$ws.Range("L97").Value = [double]"47494.35"
$ws.Range("M97").Value = [double]"40549.222534585002"
$ws.Range("N97").Value = [double]"122494.49264665099"
$ws.Range("O97").Value = [double]"14912.2101420503"
$ws.Range("P97").Value = [double]"36714.757949057297"
$ws.Range("C98").Value = [double]"0.95741497876546999"
$ws.Range("D98").Value = [double]"0.387392297941"
$ws.Range("E98").Value = [double]"1.22358417503261"
$ws.Range("F98").Value = [double]"7.0028450885262501"
$ws.Range("G98").Value = [double]"9.5712365402653408"
$ws.Range("H98").Value = [double]"0.96896921420226001"
$ws.Range("I98").Value = [double]"3.2033965116730498"
$ws.Range("J98").Value = [double]"16910.5066277642"
$ws.Range("L98").Value = [double]"36877.980000000003"
$ws.Range("M98").Value = [double]"89883.419029881203"
$ws.Range("N98").Value = [double]"147984.130643704"
$ws.Range("O98").Value = [double]"11334.4942674058"
$ws.Range("P98").Value = [double]"34044.705933990597"

